$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Extend the two "── ... ──" console banner separators with more dashes.
# ---------------------------------------------------------------------
$old1 = "## ── Attaching packages ─────────────────────────────────────────────────────────── tidyverse 1.3.0 ──"
$new1 = "## ── Attaching packages ──────────────────────────────────────────────────────────────────────────────────────── tidyverse 1.3.0 ──"
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

$old2 = "## ── Conflicts ────────────────────────────────────────────────────────────── tidyverse_conflicts() ──"
$new2 = "## ── Conflicts ─────────────────────────────────────────────────────────────────────────────────────────── tidyverse_conflicts() ──"
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Rewrite the "A method to avoid splitting a single genome ..." paragraph
#    (replicate the threshold to get one OTU per genome/species 100 times:
#    XX% placeholders become concrete bolded values) and bold the four
#    threshold numbers (5.5, 2.5, 4.0, 3.5).
# ---------------------------------------------------------------------
$oldPara = "A method to avoid splitting a single genome into multiple units of inference is to cluster 16S rRNA gene sequences together that are similar. However, this also increases the risk of lumping together genes from different species that are similar to each other. Therefore, I assessed the impact of the threshold used to define clusters of 16S rRNA genes on the propopensity to lump species together and split genome apart. I identified the threshold where 90% of bacterial species would be represented by a single OTU. For full length 16S rRNA gene sequences, I found that at a threshold of XX%, 90% of the species would be represented by a single OTU. Similarly, thresholds of XX, XX, and XX% were observed for the V4, V3-V4, and V4-V5 regions. However, at these thresholds, multiple species could be represented by the same OTU. At the highest level of resolution, XX% of the species shared a 16S rRNA gene sequence variant with another species. Given the risk of splitting a genome into multiple OTUs is more biologically problematic than lumping species together, larger thresholds are advisable."
$newPara = "A method to avoid splitting a single genome into multiple units of inference is to cluster 16S rRNA gene sequences together that are similar. However, this also increases the risk of lumping together genes from different species that are similar to each other. Therefore, I assessed the impact of the threshold used to define clusters of 16S rRNA genes on the propensity to split a genome apart or to lump species together. For full length 16S rRNA gene sequences, I found that at a threshold of 5.5%, 95% of the species with 7 copies of the rrn operon would be represented by a single OTU. Similarly, thresholds of 2.5, 4.0, and 3.5% were observed for the V4, V3-V4, and V4-V5 regions, respectively. However, at these thresholds, multiple species could be represented by the same OTU. At the highest level of resolution, XX% of the species shared a 16S rRNA gene sequence variant with another species. Given the risk of splitting a genome into multiple OTUs is more biologically problematic than lumping species together, larger thresholds are advisable."

# Locate the target paragraph; keep the Paragraph object so its .Range
# reflects subsequent in-place edits without having to re-scan the document.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("A method to avoid splitting a single genome")) {
        $targetPara = $p
        break
    }
}

$targetPara.Range.Find.Execute($oldPara, $false, $false, $false, $false, $false, $true, 1, $false, $newPara, 2) | Out-Null

# Bold the four threshold values that were substituted in for "XX".
foreach ($boldText in "5.5", "2.5", "4.0", "3.5") {
    $scopedRange = $targetPara.Range
    $found = $scopedRange.Find.Execute($boldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $scopedRange.Font.Bold = $true
    }
}
